$d = $word.ActiveDocument

# --- 1. The _GoBack bookmark currently sits at the end of "It is easy to
#        gain permission from the authorities."; it is about to be
#        recreated on the freshly inserted paragraph below, so drop the
#        old one first (its name would otherwise be ambiguous). ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- 2. Rebuild the title paragraph (add rFonts/enlarge the paragraph-mark
#        run properties) and insert the new introductory paragraph right
#        after it, carrying the _GoBack bookmark with it. ---
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range

$introText = "This document is used to analyze all five Porter" + [char]0x2019 + "s forces to explain the connection between the new entrants and how their represent a threat for our business. This document also explain customers" + [char]0x2019 + " power and what a substitute produce is."

$xml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:jc w:val="center"/>
<w:rPr><w:rFonts w:cs="Tahoma"/><w:b/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Five Forces for </w:t></w:r>
<w:r><w:rPr><w:rFonts w:cs="Tahoma"/><w:b/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>Zebras Burn Pub</w:t></w:r>
</w:p>
<w:p>
<w:pPr>
<w:jc w:val="both"/>
<w:rPr><w:szCs w:val="24"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:rFonts w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{0}</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@ -f $introText

$titleRange.InsertXML($xml)
